$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_28_2_0"
$ws.Range("B2").Value = 0.9999549481215966
$ws.Range("C2").Value = 0.999063942309954
$ws.Range("D2").Value = 0.9999609795966211
$ws.Range("E2").Value = 0.9997635270089874
$ws.Range("F2").Value = 0.9998535284158813
$ws.Range("G2").Value = 0.00004205395094067856
$ws.Range("H2").Value = 0.000873768765474148
$ws.Range("I2").Value = 0.00003055746665345906
$ws.Range("J2").Value = 0.0002322797236735885
$ws.Range("K2").Value = 0.0001314185951635238
$ws.Range("L2").Value = 0.000423683594587255
$ws.Range("M2").Value = 0.006484901768005322
$ws.Range("N2").Value = 1.000051487861032
$ws.Range("O2").Value = 0.00676097746186659
$ws.Range("P2").Value = 110.1531144355958
$ws.Range("Q2").Value = 165.0025265546649

$ws.Range("A3").Value = "model_28_2_1"
$ws.Range("B3").Value = 0.9999549431951168
$ws.Range("C3").Value = 0.9990639283859557
$ws.Range("D3").Value = 0.9999609768448401
$ws.Range("E3").Value = 0.9997635037450215
$ws.Range("F3").Value = 0.9998535123979638
$ws.Range("G3").Value = 0.00004205854959330355
$ws.Range("H3").Value = 0.0008737817629153589
$ws.Range("I3").Value = 0.00003055962161469406
$ws.Range("J3").Value = 0.0002323025751102925
$ws.Range("K3").Value = 0.0001314329669081595
$ws.Range("L3").Value = 0.0004236865849936197
$ws.Range("M3").Value = 0.006485256324410282
$ws.Range("N3").Value = 1.000051493491295
$ws.Range("O3").Value = 0.006761347112471741
$ws.Range("P3").Value = 110.1528957450278
$ws.Range("Q3").Value = 165.0023078640968

$ws.Range("A4").Value = "model_28_2_22"
$ws.Range("B4").Value = 0.9999549410293126
$ws.Range("C4").Value = 0.9990639202211722
$ws.Range("D4").Value = 0.9999609713408426
$ws.Range("E4").Value = 0.9997634914641947
$ws.Range("F4").Value = 0.9998535053561581
$ws.Range("G4").Value = 0.00004206057127643839
$ws.Range("H4").Value = 0.0008737893843822983
$ws.Range("I4").Value = 0.00003056393187822458
$ws.Range("J4").Value = 0.0002323146381667051
$ws.Range("K4").Value = 0.0001314392850224648
$ws.Range("L4").Value = 0.0004236862227404525
$ws.Range("M4").Value = 0.006485412190172526
$ws.Range("N4").Value = 1.0000514959665
$ws.Range("O4").Value = 0.00676150961376217
$ws.Range("P4").Value = 110.1527996107309
$ws.Range("Q4").Value = 165.0022117297999

$ws.Range("A5").Value = "model_28_2_21"
$ws.Range("B5").Value = 0.9999549410293126
$ws.Range("C5").Value = 0.9990639202211722
$ws.Range("D5").Value = 0.9999609713408426
$ws.Range("E5").Value = 0.9997634914641947
$ws.Range("F5").Value = 0.9998535053561581
$ws.Range("G5").Value = 0.00004206057127643839
$ws.Range("H5").Value = 0.0008737893843822983
$ws.Range("I5").Value = 0.00003056393187822458
$ws.Range("J5").Value = 0.0002323146381667051
$ws.Range("K5").Value = 0.0001314392850224648
$ws.Range("L5").Value = 0.0004236862227404525
$ws.Range("M5").Value = 0.006485412190172526
$ws.Range("N5").Value = 1.0000514959665
$ws.Range("O5").Value = 0.00676150961376217
$ws.Range("P5").Value = 110.1527996107309
$ws.Range("Q5").Value = 165.0022117297999

$ws.Range("A6").Value = "model_28_2_20"
$ws.Range("B6").Value = 0.9999549410293126
$ws.Range("C6").Value = 0.9990639202211722
$ws.Range("D6").Value = 0.9999609713408426
$ws.Range("E6").Value = 0.9997634914641947
$ws.Range("F6").Value = 0.9998535053561581
$ws.Range("G6").Value = 0.00004206057127643839
$ws.Range("H6").Value = 0.0008737893843822983
$ws.Range("I6").Value = 0.00003056393187822458
$ws.Range("J6").Value = 0.0002323146381667051
$ws.Range("K6").Value = 0.0001314392850224648
$ws.Range("L6").Value = 0.0004236862227404525
$ws.Range("M6").Value = 0.006485412190172526
$ws.Range("N6").Value = 1.0000514959665
$ws.Range("O6").Value = 0.00676150961376217
$ws.Range("P6").Value = 110.1527996107309
$ws.Range("Q6").Value = 165.0022117297999

$ws.Range("A7").Value = "model_28_2_19"
$ws.Range("B7").Value = 0.9999549410293126
$ws.Range("C7").Value = 0.9990639202211722
$ws.Range("D7").Value = 0.9999609713408426
$ws.Range("E7").Value = 0.9997634914641947
$ws.Range("F7").Value = 0.9998535053561581
$ws.Range("G7").Value = 0.00004206057127643839
$ws.Range("H7").Value = 0.0008737893843822983
$ws.Range("I7").Value = 0.00003056393187822458
$ws.Range("J7").Value = 0.0002323146381667051
$ws.Range("K7").Value = 0.0001314392850224648
$ws.Range("L7").Value = 0.0004236862227404525
$ws.Range("M7").Value = 0.006485412190172526
$ws.Range("N7").Value = 1.0000514959665
$ws.Range("O7").Value = 0.00676150961376217
$ws.Range("P7").Value = 110.1527996107309
$ws.Range("Q7").Value = 165.0022117297999

$ws.Range("A8").Value = "model_28_2_18"
$ws.Range("B8").Value = 0.9999549410293126
$ws.Range("C8").Value = 0.9990639202211722
$ws.Range("D8").Value = 0.9999609713408426
$ws.Range("E8").Value = 0.9997634914641947
$ws.Range("F8").Value = 0.9998535053561581
$ws.Range("G8").Value = 0.00004206057127643839
$ws.Range("H8").Value = 0.0008737893843822983
$ws.Range("I8").Value = 0.00003056393187822458
$ws.Range("J8").Value = 0.0002323146381667051
$ws.Range("K8").Value = 0.0001314392850224648
$ws.Range("L8").Value = 0.0004236862227404525
$ws.Range("M8").Value = 0.006485412190172526
$ws.Range("N8").Value = 1.0000514959665
$ws.Range("O8").Value = 0.00676150961376217
$ws.Range("P8").Value = 110.1527996107309
$ws.Range("Q8").Value = 165.0022117297999

$ws.Range("A9").Value = "model_28_2_17"
$ws.Range("B9").Value = 0.9999549410293126
$ws.Range("C9").Value = 0.9990639202211722
$ws.Range("D9").Value = 0.9999609713408426
$ws.Range("E9").Value = 0.9997634914641947
$ws.Range("F9").Value = 0.9998535053561581
$ws.Range("G9").Value = 0.00004206057127643839
$ws.Range("H9").Value = 0.0008737893843822983
$ws.Range("I9").Value = 0.00003056393187822458
$ws.Range("J9").Value = 0.0002323146381667051
$ws.Range("K9").Value = 0.0001314392850224648
$ws.Range("L9").Value = 0.0004236862227404525
$ws.Range("M9").Value = 0.006485412190172526
$ws.Range("N9").Value = 1.0000514959665
$ws.Range("O9").Value = 0.00676150961376217
$ws.Range("P9").Value = 110.1527996107309
$ws.Range("Q9").Value = 165.0022117297999

$ws.Range("A10").Value = "model_28_2_16"
$ws.Range("B10").Value = 0.9999549410293126
$ws.Range("C10").Value = 0.9990639202211722
$ws.Range("D10").Value = 0.9999609713408426
$ws.Range("E10").Value = 0.9997634914641947
$ws.Range("F10").Value = 0.9998535053561581
$ws.Range("G10").Value = 0.00004206057127643839
$ws.Range("H10").Value = 0.0008737893843822983
$ws.Range("I10").Value = 0.00003056393187822458
$ws.Range("J10").Value = 0.0002323146381667051
$ws.Range("K10").Value = 0.0001314392850224648
$ws.Range("L10").Value = 0.0004236862227404525
$ws.Range("M10").Value = 0.006485412190172526
$ws.Range("N10").Value = 1.0000514959665
$ws.Range("O10").Value = 0.00676150961376217
$ws.Range("P10").Value = 110.1527996107309
$ws.Range("Q10").Value = 165.0022117297999

$ws.Range("A11").Value = "model_28_2_15"
$ws.Range("B11").Value = 0.9999549410293126
$ws.Range("C11").Value = 0.9990639202211722
$ws.Range("D11").Value = 0.9999609713408426
$ws.Range("E11").Value = 0.9997634914641947
$ws.Range("F11").Value = 0.9998535053561581
$ws.Range("G11").Value = 0.00004206057127643839
$ws.Range("H11").Value = 0.0008737893843822983
$ws.Range("I11").Value = 0.00003056393187822458
$ws.Range("J11").Value = 0.0002323146381667051
$ws.Range("K11").Value = 0.0001314392850224648
$ws.Range("L11").Value = 0.0004236862227404525
$ws.Range("M11").Value = 0.006485412190172526
$ws.Range("N11").Value = 1.0000514959665
$ws.Range("O11").Value = 0.00676150961376217
$ws.Range("P11").Value = 110.1527996107309
$ws.Range("Q11").Value = 165.0022117297999

$ws.Range("A12").Value = "model_28_2_14"
$ws.Range("B12").Value = 0.9999549410293126
$ws.Range("C12").Value = 0.9990639202211722
$ws.Range("D12").Value = 0.9999609713408426
$ws.Range("E12").Value = 0.9997634914641947
$ws.Range("F12").Value = 0.9998535053561581
$ws.Range("G12").Value = 0.00004206057127643839
$ws.Range("H12").Value = 0.0008737893843822983
$ws.Range("I12").Value = 0.00003056393187822458
$ws.Range("J12").Value = 0.0002323146381667051
$ws.Range("K12").Value = 0.0001314392850224648
$ws.Range("L12").Value = 0.0004236862227404525
$ws.Range("M12").Value = 0.006485412190172526
$ws.Range("N12").Value = 1.0000514959665
$ws.Range("O12").Value = 0.00676150961376217
$ws.Range("P12").Value = 110.1527996107309
$ws.Range("Q12").Value = 165.0022117297999

$ws.Range("A13").Value = "model_28_2_13"
$ws.Range("B13").Value = 0.9999549410293126
$ws.Range("C13").Value = 0.9990639202211722
$ws.Range("D13").Value = 0.9999609713408426
$ws.Range("E13").Value = 0.9997634914641947
$ws.Range("F13").Value = 0.9998535053561581
$ws.Range("G13").Value = 0.00004206057127643839
$ws.Range("H13").Value = 0.0008737893843822983
$ws.Range("I13").Value = 0.00003056393187822458
$ws.Range("J13").Value = 0.0002323146381667051
$ws.Range("K13").Value = 0.0001314392850224648
$ws.Range("L13").Value = 0.0004236862227404525
$ws.Range("M13").Value = 0.006485412190172526
$ws.Range("N13").Value = 1.0000514959665
$ws.Range("O13").Value = 0.00676150961376217
$ws.Range("P13").Value = 110.1527996107309
$ws.Range("Q13").Value = 165.0022117297999

$ws.Range("A14").Value = "model_28_2_12"
$ws.Range("B14").Value = 0.9999549410293126
$ws.Range("C14").Value = 0.9990639202211722
$ws.Range("D14").Value = 0.9999609713408426
$ws.Range("E14").Value = 0.9997634914641947
$ws.Range("F14").Value = 0.9998535053561581
$ws.Range("G14").Value = 0.00004206057127643839
$ws.Range("H14").Value = 0.0008737893843822983
$ws.Range("I14").Value = 0.00003056393187822458
$ws.Range("J14").Value = 0.0002323146381667051
$ws.Range("K14").Value = 0.0001314392850224648
$ws.Range("L14").Value = 0.0004236862227404525
$ws.Range("M14").Value = 0.006485412190172526
$ws.Range("N14").Value = 1.0000514959665
$ws.Range("O14").Value = 0.00676150961376217
$ws.Range("P14").Value = 110.1527996107309
$ws.Range("Q14").Value = 165.0022117297999

$ws.Range("A15").Value = "model_28_2_11"
$ws.Range("B15").Value = 0.9999549410293126
$ws.Range("C15").Value = 0.9990639202211722
$ws.Range("D15").Value = 0.9999609713408426
$ws.Range("E15").Value = 0.9997634914641947
$ws.Range("F15").Value = 0.9998535053561581
$ws.Range("G15").Value = 0.00004206057127643839
$ws.Range("H15").Value = 0.0008737893843822983
$ws.Range("I15").Value = 0.00003056393187822458
$ws.Range("J15").Value = 0.0002323146381667051
$ws.Range("K15").Value = 0.0001314392850224648
$ws.Range("L15").Value = 0.0004236862227404525
$ws.Range("M15").Value = 0.006485412190172526
$ws.Range("N15").Value = 1.0000514959665
$ws.Range("O15").Value = 0.00676150961376217
$ws.Range("P15").Value = 110.1527996107309
$ws.Range("Q15").Value = 165.0022117297999

$ws.Range("A16").Value = "model_28_2_10"
$ws.Range("B16").Value = 0.9999549410293126
$ws.Range("C16").Value = 0.9990639202211722
$ws.Range("D16").Value = 0.9999609713408426
$ws.Range("E16").Value = 0.9997634914641947
$ws.Range("F16").Value = 0.9998535053561581
$ws.Range("G16").Value = 0.00004206057127643839
$ws.Range("H16").Value = 0.0008737893843822983
$ws.Range("I16").Value = 0.00003056393187822458
$ws.Range("J16").Value = 0.0002323146381667051
$ws.Range("K16").Value = 0.0001314392850224648
$ws.Range("L16").Value = 0.0004236862227404525
$ws.Range("M16").Value = 0.006485412190172526
$ws.Range("N16").Value = 1.0000514959665
$ws.Range("O16").Value = 0.00676150961376217
$ws.Range("P16").Value = 110.1527996107309
$ws.Range("Q16").Value = 165.0022117297999

$ws.Range("A17").Value = "model_28_2_9"
$ws.Range("B17").Value = 0.9999549410293126
$ws.Range("C17").Value = 0.9990639202211722
$ws.Range("D17").Value = 0.9999609713408426
$ws.Range("E17").Value = 0.9997634914641947
$ws.Range("F17").Value = 0.9998535053561581
$ws.Range("G17").Value = 0.00004206057127643839
$ws.Range("H17").Value = 0.0008737893843822983
$ws.Range("I17").Value = 0.00003056393187822458
$ws.Range("J17").Value = 0.0002323146381667051
$ws.Range("K17").Value = 0.0001314392850224648
$ws.Range("L17").Value = 0.0004236862227404525
$ws.Range("M17").Value = 0.006485412190172526
$ws.Range("N17").Value = 1.0000514959665
$ws.Range("O17").Value = 0.00676150961376217
$ws.Range("P17").Value = 110.1527996107309
$ws.Range("Q17").Value = 165.0022117297999

$ws.Range("A18").Value = "model_28_2_8"
$ws.Range("B18").Value = 0.9999549410293126
$ws.Range("C18").Value = 0.9990639202211722
$ws.Range("D18").Value = 0.9999609713408426
$ws.Range("E18").Value = 0.9997634914641947
$ws.Range("F18").Value = 0.9998535053561581
$ws.Range("G18").Value = 0.00004206057127643839
$ws.Range("H18").Value = 0.0008737893843822983
$ws.Range("I18").Value = 0.00003056393187822458
$ws.Range("J18").Value = 0.0002323146381667051
$ws.Range("K18").Value = 0.0001314392850224648
$ws.Range("L18").Value = 0.0004236862227404525
$ws.Range("M18").Value = 0.006485412190172526
$ws.Range("N18").Value = 1.0000514959665
$ws.Range("O18").Value = 0.00676150961376217
$ws.Range("P18").Value = 110.1527996107309
$ws.Range("Q18").Value = 165.0022117297999

$ws.Range("A19").Value = "model_28_2_7"
$ws.Range("B19").Value = 0.9999549410293126
$ws.Range("C19").Value = 0.9990639202211722
$ws.Range("D19").Value = 0.9999609713408426
$ws.Range("E19").Value = 0.9997634914641947
$ws.Range("F19").Value = 0.9998535053561581
$ws.Range("G19").Value = 0.00004206057127643839
$ws.Range("H19").Value = 0.0008737893843822983
$ws.Range("I19").Value = 0.00003056393187822458
$ws.Range("J19").Value = 0.0002323146381667051
$ws.Range("K19").Value = 0.0001314392850224648
$ws.Range("L19").Value = 0.0004236862227404525
$ws.Range("M19").Value = 0.006485412190172526
$ws.Range("N19").Value = 1.0000514959665
$ws.Range("O19").Value = 0.00676150961376217
$ws.Range("P19").Value = 110.1527996107309
$ws.Range("Q19").Value = 165.0022117297999

$ws.Range("A20").Value = "model_28_2_6"
$ws.Range("B20").Value = 0.9999549410293126
$ws.Range("C20").Value = 0.9990639202211722
$ws.Range("D20").Value = 0.9999609713408426
$ws.Range("E20").Value = 0.9997634914641947
$ws.Range("F20").Value = 0.9998535053561581
$ws.Range("G20").Value = 0.00004206057127643839
$ws.Range("H20").Value = 0.0008737893843822983
$ws.Range("I20").Value = 0.00003056393187822458
$ws.Range("J20").Value = 0.0002323146381667051
$ws.Range("K20").Value = 0.0001314392850224648
$ws.Range("L20").Value = 0.0004236862227404525
$ws.Range("M20").Value = 0.006485412190172526
$ws.Range("N20").Value = 1.0000514959665
$ws.Range("O20").Value = 0.00676150961376217
$ws.Range("P20").Value = 110.1527996107309
$ws.Range("Q20").Value = 165.0022117297999

$ws.Range("A21").Value = "model_28_2_5"
$ws.Range("B21").Value = 0.9999549410293126
$ws.Range("C21").Value = 0.9990639202211722
$ws.Range("D21").Value = 0.9999609713408426
$ws.Range("E21").Value = 0.9997634914641947
$ws.Range("F21").Value = 0.9998535053561581
$ws.Range("G21").Value = 0.00004206057127643839
$ws.Range("H21").Value = 0.0008737893843822983
$ws.Range("I21").Value = 0.00003056393187822458
$ws.Range("J21").Value = 0.0002323146381667051
$ws.Range("K21").Value = 0.0001314392850224648
$ws.Range("L21").Value = 0.0004236862227404525
$ws.Range("M21").Value = 0.006485412190172526
$ws.Range("N21").Value = 1.0000514959665
$ws.Range("O21").Value = 0.00676150961376217
$ws.Range("P21").Value = 110.1527996107309
$ws.Range("Q21").Value = 165.0022117297999

$ws.Range("A22").Value = "model_28_2_4"
$ws.Range("B22").Value = 0.9999549410293126
$ws.Range("C22").Value = 0.9990639202211722
$ws.Range("D22").Value = 0.9999609713408426
$ws.Range("E22").Value = 0.9997634914641947
$ws.Range("F22").Value = 0.9998535053561581
$ws.Range("G22").Value = 0.00004206057127643839
$ws.Range("H22").Value = 0.0008737893843822983
$ws.Range("I22").Value = 0.00003056393187822458
$ws.Range("J22").Value = 0.0002323146381667051
$ws.Range("K22").Value = 0.0001314392850224648
$ws.Range("L22").Value = 0.0004236862227404525
$ws.Range("M22").Value = 0.006485412190172526
$ws.Range("N22").Value = 1.0000514959665
$ws.Range("O22").Value = 0.00676150961376217
$ws.Range("P22").Value = 110.1527996107309
$ws.Range("Q22").Value = 165.0022117297999

$ws.Range("A23").Value = "model_28_2_3"
$ws.Range("B23").Value = 0.9999549410293126
$ws.Range("C23").Value = 0.9990639202211722
$ws.Range("D23").Value = 0.9999609713408426
$ws.Range("E23").Value = 0.9997634914641947
$ws.Range("F23").Value = 0.9998535053561581
$ws.Range("G23").Value = 0.00004206057127643839
$ws.Range("H23").Value = 0.0008737893843822983
$ws.Range("I23").Value = 0.00003056393187822458
$ws.Range("J23").Value = 0.0002323146381667051
$ws.Range("K23").Value = 0.0001314392850224648
$ws.Range("L23").Value = 0.0004236862227404525
$ws.Range("M23").Value = 0.006485412190172526
$ws.Range("N23").Value = 1.0000514959665
$ws.Range("O23").Value = 0.00676150961376217
$ws.Range("P23").Value = 110.1527996107309
$ws.Range("Q23").Value = 165.0022117297999

$ws.Range("A24").Value = "model_28_2_2"
$ws.Range("B24").Value = 0.9999549410293126
$ws.Range("C24").Value = 0.9990639202211722
$ws.Range("D24").Value = 0.9999609713408426
$ws.Range("E24").Value = 0.9997634914641947
$ws.Range("F24").Value = 0.9998535053561581
$ws.Range("G24").Value = 0.00004206057127643839
$ws.Range("H24").Value = 0.0008737893843822983
$ws.Range("I24").Value = 0.00003056393187822458
$ws.Range("J24").Value = 0.0002323146381667051
$ws.Range("K24").Value = 0.0001314392850224648
$ws.Range("L24").Value = 0.0004236862227404525
$ws.Range("M24").Value = 0.006485412190172526
$ws.Range("N24").Value = 1.0000514959665
$ws.Range("O24").Value = 0.00676150961376217
$ws.Range("P24").Value = 110.1527996107309
$ws.Range("Q24").Value = 165.0022117297999

$ws.Range("A25").Value = "model_28_2_23"
$ws.Range("B25").Value = 0.9999549410293126
$ws.Range("C25").Value = 0.9990639202211722
$ws.Range("D25").Value = 0.9999609713408426
$ws.Range("E25").Value = 0.9997634914641947
$ws.Range("F25").Value = 0.9998535053561581
$ws.Range("G25").Value = 0.00004206057127643839
$ws.Range("H25").Value = 0.0008737893843822983
$ws.Range("I25").Value = 0.00003056393187822458
$ws.Range("J25").Value = 0.0002323146381667051
$ws.Range("K25").Value = 0.0001314392850224648
$ws.Range("L25").Value = 0.0004236862227404525
$ws.Range("M25").Value = 0.006485412190172526
$ws.Range("N25").Value = 1.0000514959665
$ws.Range("O25").Value = 0.00676150961376217
$ws.Range("P25").Value = 110.1527996107309
$ws.Range("Q25").Value = 165.0022117297999

$ws.Range("A26").Value = "model_28_2_24"
$ws.Range("B26").Value = 0.9999549410293126
$ws.Range("C26").Value = 0.9990639202211722
$ws.Range("D26").Value = 0.9999609713408426
$ws.Range("E26").Value = 0.9997634914641947
$ws.Range("F26").Value = 0.9998535053561581
$ws.Range("G26").Value = 0.00004206057127643839
$ws.Range("H26").Value = 0.0008737893843822983
$ws.Range("I26").Value = 0.00003056393187822458
$ws.Range("J26").Value = 0.0002323146381667051
$ws.Range("K26").Value = 0.0001314392850224648
$ws.Range("L26").Value = 0.0004236862227404525
$ws.Range("M26").Value = 0.006485412190172526
$ws.Range("N26").Value = 1.0000514959665
$ws.Range("O26").Value = 0.00676150961376217
$ws.Range("P26").Value = 110.1527996107309
$ws.Range("Q26").Value = 165.0022117297999
